$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.331.69'
$ws.Range("E2").Value = '  +3.48%  '
$ws.Range("D3").Value = '3.624.33'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.10'
$ws.Range("E5").Value = '  +4.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '659.54'
$ws.Range("E6").Value = '  +2.01%  '
$ws.Range("E7").Value = '  +20.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.420'
$ws.Range("E8").Value = '  +7.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("E9").Value = '  +11.11%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").Value = '3.623.76'
$ws.Range("E11").Value = '  +2.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.32'
$ws.Range("E12").Value = '  +6.29%  '
$ws.Range("E13").Value = '  +2.78%  '
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '4.299.38'
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = '98.266.89'
$ws.Range("E16").Value = '  +3.99%  '
$ws.Range("E17").Value = '  +4.55%  '
$ws.Range("D18").Value = '3.621.95'
$ws.Range("E18").Value = '  +2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.14'
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.86'
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.13'
$ws.Range("E21").Value = '  +3.45%  '
$ws.Range("E22").Value = '  +14.25%  '
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '516.06'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.94'
$ws.Range("E26").Value = '  +7.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '100.17'
$ws.Range("E27").Value = '  +6.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.05'
$ws.Range("E28").Value = '  +6.48%  '
$ws.Range("D29").Value = '3.820.93'
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("E30").Value = '  +13.43%  '
$ws.Range("E31").Value = '  +2.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.85'
$ws.Range("E32").Value = '  +6.67%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.186'
$ws.Range("E34").Value = '  +5.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.989'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.99'
$ws.Range("E36").Value = '  +1.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.93'
$ws.Range("E37").Value = '  +9.62%  '
$ws.Range("E38").Value = '  +4.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '613.19'
$ws.Range("E39").Value = '  +11.05%  '
$ws.Range("E40").Value = '  +9.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.01'
$ws.Range("E41").Value = '  +15.70%  '
$ws.Range("E42").Value = '  +3.77%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.929'
$ws.Range("E44").Value = '  +4.38%  '
$ws.Range("E45").Value = '  +8.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0442'
$ws.Range("E46").Value = '  +8.97%  '
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.69'
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.62'
$ws.Range("E49").Value = '  +8.63%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.406'
$ws.Range("E50").Value = '  +40.34%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '32.99'
$ws.Range("E51").Value = '  -2.91%  '
